$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 398
$ws.Range("F5").Value = 1349
$ws.Range("F6").Value = 242
$ws.Range("F7").Value = 2582
$ws.Range("F8").Value = 954
$ws.Range("F9").Value = 18978
$ws.Range("F10").Value = 61
$ws.Range("F11").Value = 2038
$ws.Range("F12").Value = 693
$ws.Range("F14").Value = 362
$ws.Range("F15").Value = 625
$ws.Range("F16").Value = 205
$ws.Range("F17").Value = 216
$ws.Range("F19").Value = 331
$ws.Range("F20").Value = 53
$ws.Range("F21").Value = 221
$ws.Range("F23").Value = 132
$ws.Range("F24").Value = 8

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 178
$ws.Range("F9").Value = 115
$ws.Range("F10").Value = 243
$ws.Range("F11").Value = 243
$ws.Range("F19").Value = 32

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5929
$ws.Range("F3").Value = 597
$ws.Range("F4").Value = 566

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5929
$ws.Range("F4").Value = 597
$ws.Range("F5").Value = 566
$ws.Range("F9").Value = 1349
$ws.Range("F11").Value = 242
$ws.Range("F12").Value = 178
$ws.Range("F14").Value = 2582
$ws.Range("F15").Value = 954
$ws.Range("F16").Value = 18978
$ws.Range("F19").Value = 61
$ws.Range("F20").Value = 115
$ws.Range("F21").Value = 243
$ws.Range("F22").Value = 243
$ws.Range("F23").Value = 2039
$ws.Range("F24").Value = 693
$ws.Range("F26").Value = 362
$ws.Range("F27").Value = 625
$ws.Range("F28").Value = 205
$ws.Range("F29").Value = 216
$ws.Range("F33").Value = 331
$ws.Range("F34").Value = 53
$ws.Range("F36").Value = 221
$ws.Range("F39").Value = 132
$ws.Range("F40").Value = 32
$ws.Range("F42").Value = 8

